$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$shape = $s.Shapes.Item(6)
$shape.TextFrame.TextRange.Runs(1).Text = "ApproveExpense"
